$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) - F column updates ("想去人数" counts)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 3228
$ws1.Range("F5").Value = 2285
$ws1.Range("F8").Value = 1202
$ws1.Range("F9").Value = 1043
$ws1.Range("F10").Value = 260
$ws1.Range("F16").Value = 8072
$ws1.Range("F17").Value = 354
$ws1.Range("F25").Value = 1138
$ws1.Range("F27").Value = 1866
$ws1.Range("F28").Value = 532
$ws1.Range("F30").Value = 1687
$ws1.Range("F34").Value = 6
$ws1.Range("F35").Value = 56
$ws1.Range("F37").Value = 289
$ws1.Range("F40").Value = 361

# Sheet "演出" (sheet2) - F column update
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F6").Value = 1

# Sheet "全部类型" (sheet4) - F column updates (aggregated view of all sheets)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 3228
$ws4.Range("F7").Value = 2285
$ws4.Range("F10").Value = 1202
$ws4.Range("F12").Value = 1043
$ws4.Range("F13").Value = 260
$ws4.Range("F14").Value = 476
$ws4.Range("F18").Value = 8072
$ws4.Range("F19").Value = 354
$ws4.Range("F28").Value = 1138
$ws4.Range("F30").Value = 1866
$ws4.Range("F31").Value = 533
$ws4.Range("F33").Value = 1687
$ws4.Range("F37").Value = 6
$ws4.Range("F38").Value = 56
$ws4.Range("F40").Value = 289
$ws4.Range("F43").Value = 361
$ws4.Range("F44").Value = 1
